# Updated symbol list on Fri Feb 10 20:11:06 UTC 2023 with GitHub Actions
# Refreshes Price / Volume(1h) / Hora columns for the coin-ranking snapshot rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Cell="D2"; Value='308.09'},
    @{Cell="E2"; Value='-1.92%'},
    @{Cell="G2"; Value='20'},
    @{Cell="D3"; Value='40.79'},
    @{Cell="E3"; Value='-0.39%'},
    @{Cell="G3"; Value='20'},
    @{Cell="D4"; Value='5.045'},
    @{Cell="E4"; Value='-1.71%'},
    @{Cell="G4"; Value='20'},
    @{Cell="D5"; Value='0.07635'},
    @{Cell="E5"; Value='-2.88%'},
    @{Cell="G5"; Value='20'},
    @{Cell="D6"; Value='4.256'},
    @{Cell="E6"; Value='-2.13%'},
    @{Cell="G6"; Value='20'},
    @{Cell="D7"; Value='1.608'},
    @{Cell="E7"; Value='-3.03%'},
    @{Cell="G7"; Value='20'},
    @{Cell="D8"; Value='0.9082'},
    @{Cell="E8"; Value='-1.92%'},
    @{Cell="G8"; Value='20'},
    @{Cell="D9"; Value='2.400'},
    @{Cell="E9"; Value='-6.67%'},
    @{Cell="G9"; Value='20'},
    @{Cell="D10"; Value='0.1005'},
    @{Cell="E10"; Value='-6.03%'},
    @{Cell="G10"; Value='20'},
    @{Cell="D11"; Value='0.1768'},
    @{Cell="E11"; Value='-1.75%'},
    @{Cell="G11"; Value='20'},
    @{Cell="D12"; Value='0.09114'},
    @{Cell="E12"; Value='0.47%'},
    @{Cell="G12"; Value='20'},
    @{Cell="D13"; Value='0.04401'},
    @{Cell="E13"; Value='-1.55%'},
    @{Cell="G13"; Value='20'},
    @{Cell="D14"; Value='0.1052'},
    @{Cell="E14"; Value='-0.88%'},
    @{Cell="G14"; Value='20'},
    @{Cell="E15"; Value='-1.47%'},
    @{Cell="G15"; Value='20'},
    @{Cell="D16"; Value='0.005819'},
    @{Cell="E16"; Value='-1.45%'},
    @{Cell="G16"; Value='20'},
    @{Cell="D17"; Value='3.366'},
    @{Cell="E17"; Value='0.33%'},
    @{Cell="G17"; Value='20'},
    @{Cell="D18"; Value='0.3268'},
    @{Cell="E18"; Value='-2.31%'},
    @{Cell="G18"; Value='20'},
    @{Cell="D19"; Value='6.746'},
    @{Cell="E19"; Value='-6.94%'},
    @{Cell="G19"; Value='20'},
    @{Cell="D20"; Value='0.1357'},
    @{Cell="G20"; Value='20'},
    @{Cell="D21"; Value='0.2719'},
    @{Cell="E21"; Value='2.78%'},
    @{Cell="G21"; Value='20'},
    @{Cell="D22"; Value='0.04153'},
    @{Cell="E22"; Value='-0.78%'},
    @{Cell="G22"; Value='20'},
    @{Cell="D23"; Value='0.001215'},
    @{Cell="E23"; Value='-2.82%'},
    @{Cell="G23"; Value='20'},
    @{Cell="D24"; Value='0.004079'},
    @{Cell="E24"; Value='-1.88%'},
    @{Cell="G24"; Value='20'},
    @{Cell="D25"; Value='0.0001298'},
    @{Cell="E25"; Value='5.68%'},
    @{Cell="G25"; Value='20'},
    @{Cell="D26"; Value='0.0003003'},
    @{Cell="E26"; Value='-0.40%'},
    @{Cell="G26"; Value='20'},
    @{Cell="G27"; Value='20'},
    @{Cell="G28"; Value='20'},
    @{Cell="G29"; Value='20'},
    @{Cell="G30"; Value='20'},
    @{Cell="G31"; Value='20'},
    @{Cell="G32"; Value='20'},
    @{Cell="G33"; Value='20'},
    @{Cell="G34"; Value='20'},
    @{Cell="G35"; Value='20'},
    @{Cell="G36"; Value='20'},
    @{Cell="G37"; Value='20'},
    @{Cell="D38"; Value='0.02425'},
    @{Cell="E38"; Value='-1.13%'},
    @{Cell="G38"; Value='20'},
    @{Cell="D39"; Value='0.05183'},
    @{Cell="E39"; Value='-2.62%'},
    @{Cell="G39"; Value='20'},
    @{Cell="D40"; Value='0.007785'},
    @{Cell="E40"; Value='-2.70%'},
    @{Cell="G40"; Value='20'},
    @{Cell="D41"; Value='0.1310'},
    @{Cell="E41"; Value='-3.59%'},
    @{Cell="G41"; Value='20'},
    @{Cell="D42"; Value='0.007082'},
    @{Cell="E42"; Value='-7.60%'},
    @{Cell="G42"; Value='20'},
    @{Cell="D43"; Value='0.001945'},
    @{Cell="E43"; Value='2.88%'},
    @{Cell="G43"; Value='20'},
    @{Cell="D44"; Value='0.008006'},
    @{Cell="E44"; Value='-2.93%'},
    @{Cell="G44"; Value='20'},
    @{Cell="D45"; Value='0.3061'},
    @{Cell="E45"; Value='-2.00%'},
    @{Cell="G45"; Value='20'},
    @{Cell="D46"; Value='0.00006356'},
    @{Cell="E46"; Value='-6.71%'},
    @{Cell="G46"; Value='20'},
    @{Cell="D47"; Value='0.00000000749'},
    @{Cell="E47"; Value='-1.34%'},
    @{Cell="G47"; Value='20'},
    @{Cell="D48"; Value='0.004393'},
    @{Cell="E48"; Value='5.88%'},
    @{Cell="G48"; Value='20'},
    @{Cell="D49"; Value='0.005409'},
    @{Cell="E49"; Value='58.05%'},
    @{Cell="G49"; Value='20'},
    @{Cell="D50"; Value='0.00002097'},
    @{Cell="E50"; Value='-1.34%'},
    @{Cell="G50"; Value='20'},
    @{Cell="D51"; Value='0.0001997'},
    @{Cell="E51"; Value='-1.34%'},
    @{Cell="G51"; Value='20'}
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    # Force text storage so numeric-looking strings (prices, percentages,
    # the hour value) are written back exactly like the source feed - as
    # literal text, not auto-coerced into numbers/percentages.
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
}
